# Updated cryptos list with latest price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '63.950.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.64%  '
$ws.Range("D3").Value = "'" + '3.500.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.03%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'" + '586.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.45%  '
$ws.Range("D6").Value = "'" + '132.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.14%  '
$ws.Range("D7").Value = "'" + '3.500.11'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.10%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").Value = "'" + '0.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.57%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("D13").Value = "'" + '4.101.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.05%  '
$ws.Range("D14").Value = "'" + '27.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.41%  '
$ws.Range("E15").Value = '  -4.87%  '
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = "'" + '3.506.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.14%  '
$ws.Range("D18").Value = "'" + '64.076.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.53%  '
$ws.Range("D19").Value = "'" + '9.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = "'" + '14.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.04%  '
$ws.Range("E21").Value = '  -4.79%  '
$ws.Range("D22").Value = "'" + '391.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("E23").Value = '  -2.85%  '
$ws.Range("D24").Value = "'" + '3.644.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("D25").Value = "'" + '72.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -7.06%  '
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("D29").Value = "'" + '7.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.97%  '
$ws.Range("D30").Value = "'" + '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -4.31%  '
$ws.Range("D32").Value = "'" + '8.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.47%  '
$ws.Range("D33").Value = "'" + '3.506.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("D35").Value = "'" + '23.87'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("E36").Value = '  -3.72%  '
$ws.Range("D37").Value = "'" + '5.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("E38").Value = '  -5.24%  '
$ws.Range("D39").Value = "'" + '6.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("D40").Value = "'" + '166.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").Value = "'" + '0.0808'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.48%  '
$ws.Range("D42").Value = "'" + '27.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").Value = "'" + '0.813'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.98%  '
$ws.Range("D44").Value = "'" + '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = "'" + '41.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.74%  '
$ws.Range("E46").Value = '  -6.25%  '
$ws.Range("E47").Value = '  -4.57%  '
$ws.Range("D48").Value = "'" + '1.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.66%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = "'" + '2.448.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = "'" + '6.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("E51").Value = '  -2.17%  '
